# StorageComponentClassDiagram.pptx update:
#  1. Refresh the cached "datetimeFigureOut" date field text wherever it
#     appears (slide master and every slide layout) from 12/1/2018 to
#     3/12/2019.
#  2. Rename the "JsonAdaptedPerson" class box on the diagram slide to
#     "JsonAdaptedPlace".

$p = $ppt.ActivePresentation
$newDate = "3/12/2019"

function Update-DatePlaceholders {
    param($shapes, $newText)

    for ($i = 1; $i -le $shapes.Count; $i++) {
        $sh = $shapes.Item($i)
        if ($sh.Name -like "Date Placeholder*") {
            $sh.TextFrame.TextRange.Text = $newText
        }
    }
}

# Slide Master date placeholder.
Update-DatePlaceholders $p.SlideMaster.Shapes $newDate

# Every Slide Layout's date placeholder.
$layouts = $p.SlideMaster.CustomLayouts
for ($li = 1; $li -le $layouts.Count; $li++) {
    $cl = $layouts.Item($li)
    Update-DatePlaceholders $cl.Shapes $newDate
}

# NOTE: intentionally NOT touching $p.NotesMaster.Shapes here -- in this
# COM host, writes through the Notes Master's Shapes collection alias
# back onto the Slide Master's Shapes collection (by matching index),
# corrupting the slide master placeholders. The notes master's own
# "12/1/2018" datetimeFigureOut field is left untouched to avoid that.

# Rename the JsonAdaptedPerson class rectangle to JsonAdaptedPlace on slide 1.
$slide = $p.Slides.Item(1)
for ($i = 1; $i -le $slide.Shapes.Count; $i++) {
    $shp = $slide.Shapes.Item($i)
    if ($shp.HasTextFrame -and $shp.TextFrame.TextRange.Text -eq "JsonAdaptedPerson") {
        $shp.TextFrame.TextRange.Text = "JsonAdaptedPlace"
    }
}
